$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-17 Friday", "2025-01-18 Saturday"),
    @("958÷5=191, 3", "321÷8=40, 1"),
    @("128÷9=14, 2", "494÷4=123, 2"),
    @("458÷8=57, 2", "357÷6=59, 3"),
    @("943÷9=104, 7", "380÷8=47, 4"),
    @("564÷4=141, 0", "788÷3=262, 2"),
    @("717÷6=119, 3", "737÷4=184, 1"),
    @("436÷3=145, 1", "155÷2=77, 1"),
    @("525÷8=65, 5", "425÷5=85, 0"),
    @("326÷5=65, 1", "586÷2=293, 0"),
    @("239÷5=47, 4", "139÷2=69, 1"),
    @("787÷9=87, 4", "186÷4=46, 2"),
    @("199÷7=28, 3", "565÷7=80, 5"),
    @("401÷2=200, 1", "979÷3=326, 1"),
    @("538÷9=59, 7", "949÷4=237, 1"),
    @("548÷9=60, 8", "607÷4=151, 3"),
    @("355÷5=71, 0", "491÷5=98, 1"),
    @("606÷9=67, 3", "680÷3=226, 2"),
    @("706÷7=100, 6", "557÷8=69, 5"),
    @("752÷6=125, 2", "409÷7=58, 3"),
    @("850÷5=170, 0", "396÷4=99, 0"),
    @("607÷8=75, 7", "729÷7=104, 1"),
    @("458÷2=229, 0", "620÷8=77, 4"),
    @("961÷2=480, 1", "501÷3=167, 0"),
    @("662÷9=73, 5", "142÷7=20, 2"),
    @("244÷4=61, 0", "624÷6=104, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}
